# Update Financials figures on the INTU worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INTU")

# Balance Sheet section updates (most-recent period, column D)
$ws.Range("D43").Value = 235000    # Net Receivables
$ws.Range("D45").Value = 569000    # Other Current Assets
$ws.Range("D46").Value = 2422000   # Total Current Assets
$ws.Range("D48").Value = 1624000   # Property Plant and Equipment
$ws.Range("D49").Value = 1733000   # Goodwill
$ws.Range("D52").Value = 302000    # Other Assets
$ws.Range("D54").Value = 5134000   # Total Assets
$ws.Range("D59").Value = 1515000   # Other Current Liabilities
$ws.Range("D60").Value = 1743000   # Total Current Liabilities
$ws.Range("D62").Value = 500000    # Other Liabilities
$ws.Range("D66").Value = 2318000   # Total Liabilities
$ws.Range("D72").Value = 8564000   # Retained Earnings
$ws.Range("D76").Value = 2816000   # Total Stockholder Equity

# Cash Flow Statement - Capital Expenditures (row 91), all periods
$ws.Range("D91").Value = -38000
$ws.Range("E91").Value = -102000
$ws.Range("F91").Value = -416000
$ws.Range("G91").Value = -142000
$ws.Range("H91").Value = -104000
$ws.Range("I91").Value = -129000
$ws.Range("J91").Value = -135000
